# Adds new columns I (header "I0") and J (header "IF") with data for rows 2-50,
# mirroring the formatting already used for column H.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto I1:J1
# so the new header cells share the same style index as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (rows 2-50) for columns I and J ---
$data = @(
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,9),
    @(7,7),
    @(10,10),
    @(7,7),
    @(8,8),
    @(9,9),
    @(10,10),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(10,10),
    @(6,6),
    @(7,8),
    @(8,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(7,7),
    @(6,6),
    @(6,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(5,5),
    @(5,5),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(4,4)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
